# SignalAnalysis Hungarian (hu-HU) localization workbook.
# Adds one new translation-table entry ("strChkWindowPosition") to the
# sorted "Tabla13" table on sheet "hu-HU". The table is kept sorted
# alphabetically by the "Key" column (C), so the new row is inserted
# between the existing "strChkPower" (row 33) and "strDifferentiationAlgorithms"
# (row 34, shifting to row 35) entries, i.e. at worksheet row 34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tabla13")

# Insert a new worksheet row at row 34 - this shifts rows 34:192 down to
# 35:193 (formats/styles/heights move with their rows) and leaves a blank
# row 34 that inherits formatting from its neighbours.
$ws.Rows.Item(34).Insert()

# Grow the table definition (and its AutoFilter) so the new row becomes
# part of "Tabla13" again: B2:F192 -> B2:F193.
$lo.Resize($ws.Range("B2:F193"))

# Populate the new row with the new localization key/comment/value.
$ws.Range("B34").Value = "localization\strings"
$ws.Range("C34").Value = "strChkWindowPosition"
$ws.Range("D34").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("E34").Value = "Remember window position and size on startup"
$ws.Range("F34").Value = ""
